$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -10
$ws.Range("F10").Value = -8
$ws.Range("F11").Value = -5
$ws.Range("F13").Value = -4
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = -10
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = 0
$ws.Range("F26").Value = 1
